$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapeCount = $s.Shapes.Count

# ------------------------------------------------------------------
# 1) Update the text block that used to read:
#      "Interact" / "Energy" / "Pipe" / "DF_ADD1"
#    so it now reads:
#      "DF_CRX" / "PIPE"
#    (the two middle paragraphs "Pipe" and "DF_ADD1" are removed, and
#    the remaining two paragraphs get new text; trailing blank
#    paragraphs are left untouched).
#
#    The text box lives inside the top-level group shape with Id=67
#    ("群組 66"); the labelled rectangle is its 1st child (Id=68).
# ------------------------------------------------------------------
$textGroup = $null
for ($i = 1; $i -le $shapeCount; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 67) {
        $textGroup = $candidate
    }
}

$sub = $textGroup.GroupItems.Item(1)
$tr = $sub.TextFrame.TextRange

# Remove the old 3rd paragraph ("Pipe").
$para3 = $tr.Paragraphs(3,1)
$para3.Delete()

# After the removal above, the old 4th paragraph ("DF_ADD1") has
# shifted down to become the 3rd paragraph; remove it too.
$para3again = $tr.Paragraphs(3,1)
$para3again.Delete()

# Rename the two remaining text paragraphs.
$para1 = $tr.Paragraphs(1,1)
$para1.Text = "DF_CRX"
$para2 = $tr.Paragraphs(2,1)
$para2.Text = "PIPE"

# ------------------------------------------------------------------
# 2) Delete the small group shape ("群組 1", Id=2) that contained an
#    arrow connector, a line connector, and a text box with "5".
# ------------------------------------------------------------------
for ($i = 1; $i -le $shapeCount; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 2) {
        $candidate.Delete()
    }
}
